$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D19").Formula = "=C19/106*100"
$ws.Range("D20:D28").Formula = "=C20/106*100"
$ws.Range("D29").Formula = "=C29/105*100"
$ws.Range("D30:D38").Formula = "=C30/105*100"

$ws.Range("D19:D38").NumberFormat = "0.00"
